# Added user option for defining continuity
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

# New header cells for the continuity options
$ws.Range("L1").Value = "start_continuity"
$ws.Range("M1").Value = "end_continuity"

# New data cells - value "F" for each data row (2-5)
$ws.Range("L2").Value = "F"
$ws.Range("M2").Value = "F"
$ws.Range("L3").Value = "F"
$ws.Range("M3").Value = "F"
$ws.Range("L4").Value = "F"
$ws.Range("M4").Value = "F"
$ws.Range("L5").Value = "F"
$ws.Range("M5").Value = "F"

# Size the new "start_continuity" column the same way the other headers were sized
$ws.Columns.Item(12).ColumnWidth = 13.17

# Scroll the view over so column C is left-most visible, and leave the new
# L8 cell selected (matches the reviewer's on-screen state after the edit)
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$ws.Range("L8").Select()
